$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B41").Value = "SAILE"
$ws.Range("C41").Value = "2013 - 09"
$ws.Range("D41").Value = 224
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = 109724.12
$ws.Range("G41").Value = 562.2
$ws.Range("B42").Value = "SAILE"
$ws.Range("C42").Value = "2013 - 10"
$ws.Range("D42").Value = 1111
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 443737.68
$ws.Range("G42").Value = 876.61
$ws.Range("B43").Value = "SAILE"
$ws.Range("C43").Value = "2013 - 11"
$ws.Range("D43").Value = 1009
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 370738.75
$ws.Range("G43").Value = 0
$ws.Range("B44").Value = "SAILE"
$ws.Range("C44").Value = "2013 - 12"
$ws.Range("D44").Value = 2069
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 848209.07
$ws.Range("G44").Value = 0
$ws.Range("B45").Value = "SAILE"
$ws.Range("C45").Value = "2014 - 01"
$ws.Range("D45").Value = 2415
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 1314661.01
$ws.Range("G45").Value = 0
$ws.Range("B46").Value = "SAILE"
$ws.Range("C46").Value = "2014 - 02"
$ws.Range("D46").Value = 2386
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 1688685.77
$ws.Range("G46").Value = 0
$ws.Range("B47").Value = "SAILE"
$ws.Range("C47").Value = "2014 - 03"
$ws.Range("D47").Value = 2990
$ws.Range("E47").Value = 309
$ws.Range("F47").Value = 3145354.56
$ws.Range("G47").Value = 55902.06
$ws.Range("B48").Value = "SAILE"
$ws.Range("C48").Value = "2014 - 04"
$ws.Range("D48").Value = 1656
$ws.Range("E48").Value = 182
$ws.Range("F48").Value = 907141.15
$ws.Range("G48").Value = 36168.02
$ws.Range("B49").Value = "VIACONSUMO"
$ws.Range("C49").Value = "2013 - 04"
$ws.Range("D49").Value = 527
$ws.Range("E49").Value = 38
$ws.Range("F49").Value = 71421.14
$ws.Range("G49").Value = 1937.2
$ws.Range("B50").Value = "VIACONSUMO"
$ws.Range("C50").Value = "2013 - 05"
$ws.Range("D50").Value = 4559
$ws.Range("E50").Value = 916
$ws.Range("F50").Value = 605252.1497
$ws.Range("G50").Value = 57711.8424
$ws.Range("B51").Value = "VIACONSUMO"
$ws.Range("C51").Value = "2013 - 06"
$ws.Range("D51").Value = 5249
$ws.Range("E51").Value = 1786
$ws.Range("F51").Value = 550799.0347
$ws.Range("G51").Value = 120557.6858
$ws.Range("B52").Value = "VIACONSUMO"
$ws.Range("C52").Value = "2013 - 07"
$ws.Range("D52").Value = 8169
$ws.Range("E52").Value = 1429
$ws.Range("F52").Value = 788979.9274
$ws.Range("G52").Value = 100371.8541
$ws.Range("B53").Value = "VIACONSUMO"
$ws.Range("C53").Value = "2013 - 08"
$ws.Range("D53").Value = 8924
$ws.Range("E53").Value = 1329
$ws.Range("F53").Value = 980875.7532
$ws.Range("G53").Value = 88378.4368
$ws.Range("B54").Value = "VIACONSUMO"
$ws.Range("C54").Value = "2013 - 09"
$ws.Range("D54").Value = 9747
$ws.Range("E54").Value = 1397
$ws.Range("F54").Value = 1056305.0381
$ws.Range("G54").Value = 88616.0768
$ws.Range("B55").Value = "VIACONSUMO"
$ws.Range("C55").Value = "2013 - 10"
$ws.Range("D55").Value = 8991
$ws.Range("E55").Value = 1235
$ws.Range("F55").Value = 1204729.2156
$ws.Range("G55").Value = 101348.4984
$ws.Range("B56").Value = "VIACONSUMO"
$ws.Range("C56").Value = "2013 - 11"
$ws.Range("D56").Value = 8916
$ws.Range("E56").Value = 1165
$ws.Range("F56").Value = 1278475.8359
$ws.Range("G56").Value = 77048.5509
$ws.Range("B57").Value = "VIACONSUMO"
$ws.Range("C57").Value = "2013 - 12"
$ws.Range("D57").Value = 8317
$ws.Range("E57").Value = 1019
$ws.Range("F57").Value = 1138887.6118
$ws.Range("G57").Value = 85378.5286
$ws.Range("B58").Value = "VIACONSUMO"
$ws.Range("C58").Value = "2014 - 01"
$ws.Range("D58").Value = 10543
$ws.Range("E58").Value = 1315
$ws.Range("F58").Value = 1355577.7394
$ws.Range("G58").Value = 98554.6606
$ws.Range("B59").Value = "VIACONSUMO"
$ws.Range("C59").Value = "2014 - 02"
$ws.Range("D59").Value = 8696
$ws.Range("E59").Value = 1135
$ws.Range("F59").Value = 1038415.0607
$ws.Range("G59").Value = 79949.1056
$ws.Range("B60").Value = "VIACONSUMO"
$ws.Range("C60").Value = "2014 - 03"
$ws.Range("D60").Value = 10628
$ws.Range("E60").Value = 1516
$ws.Range("F60").Value = 1284769.434
$ws.Range("G60").Value = 101190.566
$ws.Range("B61").Value = "VIACONSUMO"
$ws.Range("C61").Value = "2014 - 04"
$ws.Range("D61").Value = 7240
$ws.Range("E61").Value = 1077
$ws.Range("F61").Value = 880941.6883
$ws.Range("G61").Value = 69881.2052

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:G61"))
